$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NSW holidays list gained two new "Anzac Day (additional)" observance rows:
#   - 2026-04-27 (Anzac Day 2026-04-25 falls on a Saturday -> observed Monday)
#   - 2027-04-26 (Anzac Day 2027-04-25 falls on a Sunday   -> observed Monday)
# Insert the rows (shifting everything below down) and fill in the two cells,
# taking care that the date-like text in column A is not auto-converted to a
# real Excel date serial number.

function Set-HolidayRow([int]$row, [string]$dateText, [string]$name) {
    $ws.Rows("$row`:$row").Insert()

    $dateCell = $ws.Range("A$row")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dateText
    $dateCell.ClearFormats()

    $ws.Range("B$row").Value = $name
}

# Row 9: new "Anzac Day (additional)" entry for 2026, pushing the existing
# 2026-06-08 King's Birthday row (and everything after it) down by one.
Set-HolidayRow 9 "2026-04-27" "Anzac Day (additional)"

# After the insert above, the 2027-04-25 Anzac Day row (originally row 21)
# is now row 22. Insert the matching 2027 additional-observance row right
# after it, at row 23.
Set-HolidayRow 23 "2027-04-26" "Anzac Day (additional)"
